$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update firstName data rows (leading apostrophe preserves the existing
# quote-prefix cell style already applied to A2:A4)
$ws.Range("A2").Value = "'Lekha"
$ws.Range("A3").Value = "'Priya"
$ws.Range("A4").Value = "'Hari"

# Remove old row 5 data (Joel) since now only 4 rows of data remain
$ws.Range("A5").ClearContents()

# Delete entire column B (removes lastName column/data)
$ws.Range("B:B").Delete()

# Update selection to match target (whole column B selected)
$ws.Range("B1:B1048576").Select()
